$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the temperature labels in column C (shared strings get rewritten by Excel on save).
$ws.Range("C3").Value = "4 ˚C"
$ws.Range("C4").Value = "7.5 ˚C"
$ws.Range("C5").Value = "10 ˚C"

# Update the saved view state for the sheet (scrolled down and a new selection).
$ws.Range("C6").Select()
